# The "Requisitos" (prerequisites) list on this sheet had two entries, in rows 23 and 24:
#   Row 23: LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)
#   Row 24: LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)
#
# The source content was reordered so the LOT2038 entry now comes first, followed by the
# LOT2028 entry. Swap the two lines (columns B and C both carry the same text) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lot2038 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"
$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

$ws.Range("B23").Value = $lot2038
$ws.Range("C23").Value = $lot2038

$ws.Range("B24").Value = $lot2028
$ws.Range("C24").Value = $lot2028
